$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = "../../NFDataCfg/Ini/Scene/1.xml"
$ws.Range("B3").Value = "../../NFDataCfg/Ini/Scene/2.xml"
$ws.Range("B4").Value = "../../NFDataCfg/Ini/Scene/3.xml"
$ws.Range("B5").Value = "../../NFDataCfg/Ini/Scene/4.xml"
$ws.Range("B6").Value = "../../NFDataCfg/Ini/Scene/5.xml"
$ws.Range("B7").Value = "../../NFDataCfg/Ini/Scene/6.xml"
